# Change date of first Assembly lab: 9/6/22 -> 9/13/22.
#
# In the real edit, the "6" inside the original single run
# <w:r><w:t>9/6/22</w:t></w:r> was selected and retyped as "13", which
# splits the text into three runs: "9/", "13", "/22". We reproduce that
# exact run structure with Range.InsertXML rather than a plain
# Find/Replace (which would just collapse everything back into one run).

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("9/6/22", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:r><w:t>9/</w:t></w:r>' + `
           '<w:r><w:t>13</w:t></w:r>' + `
           '<w:r><w:t>/22</w:t></w:r>' + `
           '</w:p>'
    $rng.InsertXML($xml)
}
